$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.016000032424926758
$ws.Range("D3").Value = 0.01699995994567871
$ws.Range("D4").Value = 0.016000032424926758
$ws.Range("D5").Value = 0.016000032424926758
$ws.Range("D6").Value = 0.017999887466430664
$ws.Range("D7").Value = 0.015000104904174805
$ws.Range("D8").Value = 0.01699995994567871
$ws.Range("D9").Value = 0.0010001659393310547
$ws.Range("D10").Value = 0.018999814987182617
$ws.Range("D11").Value = 0.015999794006347656
